$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row describing "Finding Pairs with a Certain sum"
$ws.Range("A34").Value = 1865
$ws.Range("B34").Value = "Finding Pairs with a Certain sum"

# Update the "Type" column for "Find Lucky Integer in an Array" (row 31)
$ws.Range("C31").Value = "Frequency Table/Dictionary/HashMap"

$ws.Range("C34").Value = "HashMap/Dictionary/Frequency Table"
$ws.Range("D34").Value = "2 int arrays and a Dictionary<int,int> frequency table, update the frequency table in Add(), search the complement key in Count"
$ws.Range("E34").Value = "C# Remove() = CPP erase()"

# Keep selection/view consistent with the new last row
$ws.Range("E34").Select()

# Column B needs to widen to accommodate the new, longer title text
$ws.Columns.Item(2).ColumnWidth = 36
